# Auto-generated edit script: apply numeric updates from the commit diff
# to the Hades_Profits workbook (sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

# ALC!17 (G=38956)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 520428.9
$ws.Range("J17").Value = 520428.9
$ws.Range("L17").Value = 1561286.7
$ws.Range("N17").Value = -1561622.7

# ALC!41 (G=5478)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 483.2857
$ws.Range("I41").Value = 452
$ws.Range("J41").Value = 525
$ws.Range("K41").Value = 452
$ws.Range("L41").Value = 525
$ws.Range("M41").Value = -12
$ws.Range("N41").Value = -1405

# ALC!70 (G=12604)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3188.25
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 3751
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 11253
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -11793

# ALC!73 (G=12604)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 3188.25
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 3751
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 11253
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -13125

# ALC!82 (G=12623)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 3000
$ws.Range("J82").Value = 3000
$ws.Range("L82").Value = 9000
$ws.Range("N82").Value = -9812

# ALC!85 (G=12623)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 3000
$ws.Range("J85").Value = 3000
$ws.Range("L85").Value = 9000
$ws.Range("N85").Value = -11808

# ALC!137 (G=44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5558225.5
$ws.Range("I137").Value = 9092964
$ws.Range("J137").Value = 3636.1428
$ws.Range("K137").Value = 27278892
$ws.Range("L137").Value = 10908.4284
$ws.Range("M137").Value = -27276342
$ws.Range("N137").Value = -16008.4284

# ALC!138 (G=44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 703525.1
$ws.Range("J138").Value = 1351957.1
$ws.Range("L138").Value = 4055871.3
$ws.Range("N138").Value = -4066151.3

# ARM!32 (G=44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 40055.5
$ws.Range("I32").Value = 46400.918
$ws.Range("K32").Value = 46400.918
$ws.Range("M32").Value = -46113.918

# ARM!45 (G=27714)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2539.4614
$ws.Range("I45").Value = 2677.7778
$ws.Range("J45").Value = 2228.25
$ws.Range("K45").Value = 2677.7778
$ws.Range("L45").Value = 2228.25
$ws.Range("M45").Value = -2300.7778
$ws.Range("N45").Value = -2982.25

# ARM!64 (G=10664)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# ARM!67 (G=10664)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# ARM!76 (G=10679)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 46000
$ws.Range("J76").Value = 46000
$ws.Range("L76").Value = 46000
$ws.Range("N76").Value = -46676

# ARM!79 (G=10679)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 46000
$ws.Range("J79").Value = 46000
$ws.Range("L79").Value = 46000
$ws.Range("N79").Value = -48340

# ARM!132 (G=43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 74674.97
$ws.Range("I132").Value = 64529.5
$ws.Range("J132").Value = 87161.69500000001
$ws.Range("K132").Value = 193588.5
$ws.Range("L132").Value = 261485.085
$ws.Range("M132").Value = -191058.5
$ws.Range("N132").Value = -266545.085

# CRP!31 (G=44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1711.7878
$ws.Range("I31").Value = 1302.8572
$ws.Range("J31").Value = 2427.4167
$ws.Range("K31").Value = 1302.8572
$ws.Range("L31").Value = 2427.4167
$ws.Range("M31").Value = -1007.8572
$ws.Range("N31").Value = -3017.4167

# CRP!34 (G=44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1711.7878
$ws.Range("I34").Value = 1302.8572
$ws.Range("J34").Value = 2427.4167
$ws.Range("K34").Value = 1302.8572
$ws.Range("L34").Value = 2427.4167
$ws.Range("M34").Value = -1100.8572
$ws.Range("N34").Value = -2831.4167

# CRP!58 (G=44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 30304862
$ws.Range("J58").Value = 4571.5713
$ws.Range("L58").Value = 4571.5713
$ws.Range("N58").Value = -4977.5713

# CRP!136 (G=44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 30304862
$ws.Range("J136").Value = 4571.5713
$ws.Range("L136").Value = 13714.7139
$ws.Range("N136").Value = -18814.7139

# CUL!5 (G=43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 883.65216
$ws.Range("I5").Value = 626.55554
$ws.Range("J5").Value = 1048.9286
$ws.Range("K5").Value = 1879.66662
$ws.Range("L5").Value = 3146.7858
$ws.Range("M5").Value = -1767.66662
$ws.Range("N5").Value = -3370.7858

# CUL!11 (G=4745)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1242546.2
$ws.Range("I11").Value = 1736743.6
$ws.Range("J11").Value = 501250
$ws.Range("K11").Value = 5210230.800000001
$ws.Range("L11").Value = 1503750
$ws.Range("M11").Value = -5210090.800000001
$ws.Range("N11").Value = -1504030

# CUL!118 (G=27872)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 3613.2856
$ws.Range("J118").Value = 4490
$ws.Range("L118").Value = 13470
$ws.Range("N118").Value = -15956

# CUL!122 (G=36078)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 12821233
$ws.Range("I122").Value = 484.8889
$ws.Range("J122").Value = 41667916
$ws.Range("K122").Value = 4364.0001
$ws.Range("L122").Value = 375011244
$ws.Range("M122").Value = -1914.0001
$ws.Range("N122").Value = -375016144

# CUL!135 (G=43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 883.65216
$ws.Range("I135").Value = 626.55554
$ws.Range("J135").Value = 1048.9286
$ws.Range("K135").Value = 5638.99986
$ws.Range("L135").Value = 9440.357399999999
$ws.Range("M135").Value = -3103.99986
$ws.Range("N135").Value = -14510.3574

# GSM!18 (G=4309)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# GSM!44 (G=4143)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 3333.3333
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 3333.3333
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 3333.3333
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -4525.3333

# GSM!102 (G=36169)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1392.5217
$ws.Range("I102").Value = 1418.2222
$ws.Range("J102").Value = 1300
$ws.Range("K102").Value = 1418.2222
$ws.Range("L102").Value = 1300
$ws.Range("M102").Value = 203.7778000000001
$ws.Range("N102").Value = -4544

# GSM!126 (G=36184)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1571.48
$ws.Range("I126").Value = 1276.6428
$ws.Range("J126").Value = 1946.7273
$ws.Range("K126").Value = 3829.9284
$ws.Range("L126").Value = 5840.1819
$ws.Range("M126").Value = -1359.9284
$ws.Range("N126").Value = -10780.1819

# GSM!138 (G=42325)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 45000
$ws.Range("J138").Value = 45000
$ws.Range("L138").Value = 45000
$ws.Range("N138").Value = -55280

# LTW!46 (G=5282)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 632.375
$ws.Range("I46").Value = 763.6667
$ws.Range("J46").Value = 553.6
$ws.Range("K46").Value = 763.6667
$ws.Range("M46").Value = -575.6667
$ws.Range("N46").Value = -929.6

# LTW!136 (G=44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 667326.7
$ws.Range("I136").Value = 500990
$ws.Range("J136").Value = 1000000
$ws.Range("K136").Value = 1502970
$ws.Range("L136").Value = 3000000
$ws.Range("M136").Value = -1500420
$ws.Range("N136").Value = -3005100

# WVR!126 (G=36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 841.7778
$ws.Range("I126").Value = 759.1875
$ws.Range("K126").Value = 2277.5625
$ws.Range("M126").Value = 192.4375
